# ReportingOrganisationGroup.xlsx — swap the "codeforiati:group-code" and
# "codeforiati:group-name" columns (D and E), including the header row.
# This mirrors the upstream codelist fix that re-ordered those two columns
# (group-name now appears before group-code for every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $codeCell = $ws.Cells.Item($r, 4)   # column D: codeforiati:group-code
    $nameCell = $ws.Cells.Item($r, 5)   # column E: codeforiati:group-name

    $codeVal = $codeCell.Value2
    $nameVal = $nameCell.Value2

    $codeCell.Value = $nameVal
    $nameCell.Value = $codeVal
}
